# Applies the "Vector AI" -> "Vertex AI" rename plus the accompanying
# connector re-route/resize on the Hackathon high-level design deck.
#
#   Slide 4 ("Rectangle 27"):  "Vector AI"            -> "Vertex AI container"
#   Slide 4 ("Straight Arrow Connector 31"): re-geometried + flipped
#   Slide 5 ("Rectangle 16"):  "Deploy in Vector AI"  -> "Deploy in Vertex AI"

function Get-ShapeById {
    param($Slide, [int]$Id)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $sh = $Slide.Shapes.Item($i)
        if ($sh.Id -eq $Id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# EMU -> points conversion factor used by the PowerPoint object model.
$emuPerPt = 12700

# --- Slide 4 -----------------------------------------------------------
$slide4 = $p.Slides.Item(4)

# "Vector AI" label on the small rectangle overlapping the endpoint box.
$lblShape = Get-ShapeById $slide4 28
$lblShape.TextFrame.TextRange.Text = "Vertex AI container"

# The connector used to end at shape id 2 (idx 1); it now starts there
# instead, is flipped horizontally, and is repositioned/resized.
$connShape = Get-ShapeById $slide4 32
$endpointShape = Get-ShapeById $slide4 2

# Best-effort: re-home the connection from an "end" to a "begin" attachment
# (no-op on hosts that don't support re-pointing an existing connector, but
# harmless to attempt).
try {
    $connShape.ConnectorFormat.EndDisconnect()
    $connShape.ConnectorFormat.BeginConnect($endpointShape, 1)
} catch {
}

$connShape.Left = 4298686 / $emuPerPt
$connShape.Width = 364964 / $emuPerPt
$connShape.Flip(0)   # msoFlipHorizontal -> sets flipH="1"

# --- Slide 5 -------------------------------------------------------------
$slide5 = $p.Slides.Item(5)

$deployShape = Get-ShapeById $slide5 17
$deployShape.TextFrame.TextRange.Text = "Deploy in Vertex AI"
